$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete columns E:F entirely (shifts G,H,... left by two)
$ws.Range("E:F").EntireColumn.Delete() | Out-Null

# Update the two surviving header strings (now in C4 and D4)
$ws.Range("C4").Value = "최소자리"
$ws.Range("D4").Value = "복잡성"

# Restore the active cell selection as in the target workbook
$ws.Range("O27").Select() | Out-Null
